$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2036
$ws.Range("I113").Value = 1899
$ws.Range("J113").Value = 2447
$ws.Range("K113").Value = 1899
$ws.Range("L113").Value = 2447
$ws.Range("M113").Value = 1355
$ws.Range("N113").Value = -8955

$ws.Range("H138").Value = 3444.7036
$ws.Range("I138").Value = 1965.2142
$ws.Range("K138").Value = 5895.642599999999
$ws.Range("M138").Value = -755.6425999999992

$ws.Range("H141").Value = 1685.2
$ws.Range("I141").Value = 1685.2
$ws.Range("K141").Value = 5055.6
$ws.Range("M141").Value = 124.3999999999996

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1500.58
$ws.Range("I32").Value = 1324.2391
$ws.Range("J32").Value = 3528.5
$ws.Range("K32").Value = 1324.2391
$ws.Range("L32").Value = 3528.5
$ws.Range("M32").Value = -1037.2391
$ws.Range("N32").Value = -4102.5

$ws.Range("H122").Value = 4326.9165
$ws.Range("I122").Value = 3424.818
$ws.Range("K122").Value = 10274.454
$ws.Range("M122").Value = -7824.454000000002

$ws.Range("H132").Value = 21115394
$ws.Range("I132").Value = 12504818
$ws.Range("K132").Value = 37514454
$ws.Range("M132").Value = -37511924

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2467.7856
$ws.Range("I105").Value = 2406.8
$ws.Range("J105").Value = 2620.25
$ws.Range("K105").Value = 2406.8
$ws.Range("L105").Value = 2620.25
$ws.Range("M105").Value = -659.8000000000002
$ws.Range("N105").Value = -6114.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2469.0635
$ws.Range("I31").Value = 1717.3405
$ws.Range("K31").Value = 1717.3405
$ws.Range("M31").Value = -1422.3405

$ws.Range("H34").Value = 2469.0635
$ws.Range("I34").Value = 1717.3405
$ws.Range("K34").Value = 1717.3405
$ws.Range("M34").Value = -1515.3405

$ws.Range("H35").Value = 2413.5715
$ws.Range("I35").Value = 1131.6666
$ws.Range("K35").Value = 1131.6666
$ws.Range("M35").Value = -837.6666

$ws.Range("H45").Value = 8999.25
$ws.Range("I45").Value = 8999.25
$ws.Range("K45").Value = 8999.25
$ws.Range("M45").Value = -8406.25

$ws.Range("H58").Value = 23821886
$ws.Range("J58").Value = 1503.1428
$ws.Range("L58").Value = 1503.1428
$ws.Range("N58").Value = -1909.1428

$ws.Range("H62").Value = 2698.2
$ws.Range("I62").Value = 2245
$ws.Range("J62").Value = 3000.3333
$ws.Range("K62").Value = 2245
$ws.Range("L62").Value = 3000.3333
$ws.Range("M62").Value = -1621
$ws.Range("N62").Value = -4248.3333

$ws.Range("H65").Value = 2698.2
$ws.Range("I65").Value = 2245
$ws.Range("J65").Value = 3000.3333
$ws.Range("K65").Value = 11225
$ws.Range("L65").Value = 15001.6665
$ws.Range("M65").Value = -8105
$ws.Range("N65").Value = -21241.6665

$ws.Range("H99").Value = 3480
$ws.Range("I99").Value = 3166.6667
$ws.Range("J99").Value = 3614.2856
$ws.Range("K99").Value = 3166.6667
$ws.Range("L99").Value = 3614.2856
$ws.Range("M99").Value = -1668.6667
$ws.Range("N99").Value = -6610.2856

$ws.Range("H105").Value = 10205282
$ws.Range("I105").Value = 10205282
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 10205282
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -10203535
$ws.Range("N105").ClearContents()

$ws.Range("H126").Value = 3480
$ws.Range("I126").Value = 3166.6667
$ws.Range("J126").Value = 3614.2856
$ws.Range("K126").Value = 9500.000100000001
$ws.Range("L126").Value = 10842.8568
$ws.Range("M126").Value = -7030.000100000001
$ws.Range("N126").Value = -15782.8568

$ws.Range("H134").Value = 10871331
$ws.Range("I134").Value = 11906342
$ws.Range("J134").Value = 3724
$ws.Range("K134").Value = 35719026
$ws.Range("L134").Value = 11172
$ws.Range("M134").Value = -35716491
$ws.Range("N134").Value = -16242

$ws.Range("H136").Value = 23821886
$ws.Range("J136").Value = 1503.1428
$ws.Range("L136").Value = 4509.428400000001
$ws.Range("N136").Value = -9609.428400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4601463
$ws.Range("I4").Value = 6389639
$ws.Range("J4").Value = 3295.7144
$ws.Range("K4").Value = 19168917
$ws.Range("L4").Value = 9887.143199999999
$ws.Range("M4").Value = -19168805
$ws.Range("N4").Value = -10111.1432

$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 4543.643
$ws.Range("I102").Value = 2885.4614
$ws.Range("K102").Value = 2885.4614
$ws.Range("M102").Value = -1263.4614

$ws.Range("H113").Value = 27643.83
$ws.Range("I113").Value = 31571.5
$ws.Range("K113").Value = 31571.5
$ws.Range("M113").Value = -29401.5

$ws.Range("H132").Value = 3133667
$ws.Range("I132").Value = 3580844
$ws.Range("K132").Value = 10742532
$ws.Range("M132").Value = -10740002

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 103837.69
$ws.Range("I43").Value = 9486.125
$ws.Range("K43").Value = 9486.125
$ws.Range("M43").Value = -9293.125

$ws.Range("H55").Value = 567.7646999999999
$ws.Range("J55").Value = 708.4
$ws.Range("N55").Value = -1054.4

$ws.Range("H122").Value = 2993
$ws.Range("I122").Value = 2993
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8979
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6529
$ws.Range("N122").ClearContents()

$ws.Range("H132").Value = 11911915
$ws.Range("I132").Value = 11911915
$ws.Range("K132").Value = 35735745
$ws.Range("M132").Value = -35733215

$ws.Range("H136").Value = 1945.7872
$ws.Range("I136").Value = 1853
$ws.Range("J136").Value = 2008.75
$ws.Range("K136").Value = 5559
$ws.Range("L136").Value = 6026.25
$ws.Range("M136").Value = -3009
$ws.Range("N136").Value = -11126.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H110").Value = 145499.75
$ws.Range("J110").Value = 145499.75
$ws.Range("L110").Value = 145499.75
$ws.Range("N110").Value = -153679.75

$ws.Range("H122").Value = 3557.5715
$ws.Range("I122").Value = 3380.8
$ws.Range("K122").Value = 10142.4
$ws.Range("M122").Value = -7692.400000000001

$ws.Range("H132").Value = 15156529
$ws.Range("I132").Value = 21740836
$ws.Range("J132").Value = 12620.1
$ws.Range("K132").Value = 65222508
$ws.Range("L132").Value = 37860.3
$ws.Range("M132").Value = -65219978
$ws.Range("N132").Value = -42920.3

$ws.Range("H136").Value = 17243938
$ws.Range("I136").Value = 18521120
$ws.Range("K136").Value = 55563360
$ws.Range("M136").Value = -55560810
